$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.458.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.01%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.528.71'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.35%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.91%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.72%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.36%  '

# Row 8
$ws.Range('E8').Value = '  -0.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.638'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.06%  '

# Row 10
$ws.Range('E10').Value = '  +6.48%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.75'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.65%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000282'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.86%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.29'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.73%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.093.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.57%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.524.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.31%  '

# Row 16
$ws.Range('E16').Value = '  +0.28%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.44'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.48%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.428.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.99%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.75%  '

# Row 20
$ws.Range('E20').Value = '  +1.66%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '415.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.26%  '

# Row 22
$ws.Range('E22').Value = '  +7.43%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.98%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.33%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.95%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.78%  '

# Row 27
$ws.Range('E27').Value = '  -0.87%  '

# Row 28
$ws.Range('E28').Value = '  -1.52%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.71%  '

# Row 30
$ws.Range('E30').Value = '  +1.54%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '628.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.97%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.62'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.40%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.42%  '

# Row 34
$ws.Range('E34').Value = '  +1.11%  '

# Row 35
$ws.Range('E35').Value = '  +12.97%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.76'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.66%  '

# Row 37
$ws.Range('E37').Value = '  -0.16%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.03%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.97%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.60%  '

# Row 41
$ws.Range('E41').Value = '  -0.93%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.257.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.98%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.96'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.18%  '

# Row 45
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.56'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.00%  '

# Row 46
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0423'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.58%  '

# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.40%  '

# Row 49
$ws.Range('E49').Value = '  +2.16%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.83%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '140.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.32%  '
